$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '28.268.74'
$c.ClearFormats()

$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  +2.53%  '
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '1.870.60'
$c.ClearFormats()

$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +1.29%  '
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  -0.27%  '
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '337.23'
$c.ClearFormats()

$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  +0.92%  '
$c.ClearFormats()

$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -0.37%  '
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.4709'
$c.ClearFormats()

$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +1.15%  '
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.3929'
$c.ClearFormats()

$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +1.73%  '
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '47.33'
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +2.43%  '
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.08004'
$c.ClearFormats()

$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  +0.99%  '
$c.ClearFormats()

$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  +0.67%  '
$c.ClearFormats()

$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +1.09%  '
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '1.889.05'
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  +1.97%  '
$c.ClearFormats()

$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  +0.96%  '
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '7.286'
$c.ClearFormats()

$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  +2.27%  '
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '91.25'
$c.ClearFormats()

$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  +2.49%  '
$c.ClearFormats()

$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  -0.35%  '
$c.ClearFormats()

$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +0.67%  '
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '0.06591'
$c.ClearFormats()

$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -1.06%  '
$c.ClearFormats()

$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '17.66'
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  +3.32%  '
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()

$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  -0.29%  '
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '28.278.29'
$c.ClearFormats()

$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  +2.57%  '
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '5.445'
$c.ClearFormats()

$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  +1.04%  '
$c.ClearFormats()

$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  +1.01%  '
$c.ClearFormats()

$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  -0.49%  '
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '2.110.57'
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +1.91%  '
$c.ClearFormats()

$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '159.34'
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  +0.80%  '
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '19.89'
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  +1.71%  '
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '2.152'
$c.ClearFormats()

$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +2.21%  '
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '5.500'
$c.ClearFormats()

$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  +1.60%  '
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '120.01'
$c.ClearFormats()

$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.ClearFormats()

$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  +0.27%  '
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '0.09487'
$c.ClearFormats()

$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  +0.80%  '
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '3.582'
$c.ClearFormats()

$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  -0.20%  '
$c.ClearFormats()

$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  +2.30%  '
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '5.356'
$c.ClearFormats()

$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +1.20%  '
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '0.02267'
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  +1.59%  '
$c.ClearFormats()

$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '0.06089'
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  +0.78%  '
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '8.434'
$c.ClearFormats()

$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  +1.27%  '
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '1.173'
$c.ClearFormats()

$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -0.84%  '
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '0.5965'
$c.ClearFormats()

$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  +1.15%  '
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '0.1879'
$c.ClearFormats()

$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  +0.65%  '
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '10.41'
$c.ClearFormats()

$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  +0.58%  '
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '1.303'
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +4.75%  '
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '0.5615'
$c.ClearFormats()

$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  +0.42%  '
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '12.24'
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -0.03%  '
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '1.969'
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  +3.20%  '
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '0.06893'
$c.ClearFormats()

$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  +3.02%  '
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '110.74'
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '2.021'
$c.ClearFormats()

$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  +13.26%  '
$c.ClearFormats()
